$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1024, pushing the existing data (old rows
# 1024-1046) down to 1025-1047, carrying their formatting/values along.
$ws.Rows.Item(1024).Insert()

# Populate the newly inserted row 1024 with the new weekly price entry.
$ws.Cells.Item(1024, 1).Value = 4
$ws.Cells.Item(1024, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(1024, 3).Value = "Los Lagos"
$ws.Cells.Item(1024, 4).Value = 45239
$ws.Cells.Item(1024, 5).Value = 10
$ws.Cells.Item(1024, 6).Value = 100112004
$ws.Cells.Item(1024, 7).Value = "Cebolla"
$ws.Cells.Item(1024, 8).Value = "Sin especificar"
$ws.Cells.Item(1024, 9).Value = "1a nueva(o)"
$ws.Cells.Item(1024, 10).Value = 500
$ws.Cells.Item(1024, 11).Value = 18000
$ws.Cells.Item(1024, 12).Value = 18000
$ws.Cells.Item(1024, 13).Value = 18000
$ws.Cells.Item(1024, 14).Value = "`$/malla 17 kilos"
$ws.Cells.Item(1024, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(1024, 16).Value = 1059
$ws.Cells.Item(1024, 17).Value = 17
$ws.Cells.Item(1024, 18).Value = "Hortaliza"
